# Refresh the cryptos price/volume snapshot (Price column D, Volume(1h) column E).
#
# Some Price values in column D are stored as plain-looking decimal numbers
# (e.g. "0.9990", "1.0000", "8.380") even though the source data keeps them as
# text to preserve exact digit counts/trailing zeros. Assigning such a string
# straight to .Value lets Excel silently reinterpret it as a Number and drop the
# formatting (e.g. "1.0000" -> 1, "8.380" -> 8.38), so for those specific cells we
# first switch the cell to Text number format ("@") to keep the literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textFormatCells = @('D4','D5','D6','D7','D8','D10','D11','D13','D15','D16','D20','D21','D23','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D37','D40','D41','D42','D43','D45','D46','D48','D49','D50')
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.394.66'
$ws.Range('D3').Value = '1.849.02'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D4').Value = '0.9990'
$ws.Range('D5').Value = '240.42'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '0.6297'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').Value = '1.0000'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.07623'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D10').Value = '24.54'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').Value = '0.07738'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = '1.834.65'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').Value = '5.006'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('E14').Value = '  +7.94%  '
$ws.Range('D15').Value = '0.6795'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '83.61'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('D17').Value = '2.080.95'
$ws.Range('E17').Value = '  -8.10%  '
$ws.Range('D19').Value = '29.418.32'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').Value = '228.83'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '12.45'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '7.464'
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').Value = '8.380'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').Value = '17.65'
$ws.Range('D29').Value = '1.469'
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').Value = '1.310'
$ws.Range('E30').Value = '  +4.58%  '
$ws.Range('D31').Value = '0.05646'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('D32').Value = '4.115'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('D33').Value = '4.039'
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('D34').Value = '1.852'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '1.157'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').Value = '0.7097'
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('D37').Value = '2.584'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('D38').Value = '1.235.61'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('D40').Value = '0.01798'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('D41').Value = '6.466'
$ws.Range('E41').Value = '  +4.71%  '
$ws.Range('D42').Value = '0.9075'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = '0.9997'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = '1.990.75'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = '101.47'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').Value = '66.06'
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('E47').Value = '  +4.08%  '
$ws.Range('D48').Value = '7.155'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').Value = '0.4019'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').Value = '9.032'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -0.89%  '
